$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the comment/header text in A1 (shared string content change)
$ws.Range("A1").Value = "//注释,配置npc表的士兵id"

# C1 previously had its own (duplicate) wrap-text style; re-applying the
# same formatting lets the engine re-use the existing identical style
# record instead of keeping a separate duplicate one.
$ws.Range("C1").WrapText = $true

# Move the active selection to J12, matching the saved cursor position.
$ws.Range("J12").Select() | Out-Null
